$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 18523968
$ws.Range("I111").Value = 37039824
$ws.Range("J111").Value = 8114.6665
$ws.Range("K111").Value = 111119472
$ws.Range("L111").Value = 24343.9995
$ws.Range("M111").Value = -111116405
$ws.Range("N111").Value = -30477.9995
$ws.Range("H113").Value = 45457470
$ws.Range("I113").Value = 62502000
$ws.Range("J113").Value = 5396
$ws.Range("K113").Value = 62502000
$ws.Range("L113").Value = 5396
$ws.Range("M113").Value = -62498746
$ws.Range("N113").Value = -11904
$ws.Range("H129").Value = 1071.678
$ws.Range("I129").Value = 332.54544
$ws.Range("K129").Value = 997.63632
$ws.Range("M129").Value = 4002.36368
$ws.Range("H135").Value = 69460.53
$ws.Range("I135").Value = 74064.86
$ws.Range("K135").Value = 666583.74
$ws.Range("M135").Value = -664048.74
$ws.Range("H137").Value = 1242.4
$ws.Range("I137").Value = 1086.4762
$ws.Range("J137").Value = 1606.2222
$ws.Range("K137").Value = 3259.4286
$ws.Range("L137").Value = 4818.6666
$ws.Range("M137").Value = -709.4286000000002
$ws.Range("N137").Value = -9918.6666
$ws.Range("H138").Value = 1898.8143
$ws.Range("I138").Value = 1298.0244
$ws.Range("J138").Value = 2748.2068
$ws.Range("K138").Value = 3894.0732
$ws.Range("L138").Value = 8244.6204
$ws.Range("M138").Value = 1245.9268
$ws.Range("N138").Value = -18524.6204
$ws.Range("H141").Value = 1555.7551
$ws.Range("I141").Value = 1551.76
$ws.Range("J141").Value = 1559.9166
$ws.Range("K141").Value = 4655.28
$ws.Range("L141").Value = 4679.7498
$ws.Range("M141").Value = 524.7200000000003
$ws.Range("N141").Value = -15039.7498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20968.793
$ws.Range("I32").Value = 3980.5078
$ws.Range("K32").Value = 3980.5078
$ws.Range("M32").Value = -3693.5078
$ws.Range("H44").Value = 24831.166
$ws.Range("J44").Value = 24831.166
$ws.Range("L44").Value = 24831.166
$ws.Range("N44").Value = -25807.166
$ws.Range("H55").Value = 16697.25
$ws.Range("J55").Value = 16697.25
$ws.Range("L55").Value = 16697.25
$ws.Range("N55").Value = -17327.25
$ws.Range("H61").Value = 3158.7368
$ws.Range("I61").Value = 1768.3846
$ws.Range("J61").Value = 6171.1665
$ws.Range("K61").Value = 1768.3846
$ws.Range("L61").Value = 6171.1665
$ws.Range("M61").Value = -1556.3846
$ws.Range("N61").Value = -6595.1665
$ws.Range("H74").Value = 1583
$ws.Range("I74").Value = 1106.0454
$ws.Range("J74").Value = 2536.9092
$ws.Range("K74").Value = 1106.0454
$ws.Range("L74").Value = 2536.9092
$ws.Range("M74").Value = -232.0454
$ws.Range("N74").Value = -4284.9092
$ws.Range("H77").Value = 1583
$ws.Range("I77").Value = 1106.0454
$ws.Range("J77").Value = 2536.9092
$ws.Range("K77").Value = 5530.227
$ws.Range("L77").Value = 12684.546
$ws.Range("M77").Value = -1162.227
$ws.Range("N77").Value = -21420.546
$ws.Range("H80").Value = 34626
$ws.Range("J80").Value = 34626
$ws.Range("L80").Value = 34626
$ws.Range("N80").Value = -36622
$ws.Range("H83").Value = 34626
$ws.Range("J83").Value = 34626
$ws.Range("L83").Value = 103878
$ws.Range("N83").Value = -113862
$ws.Range("H110").Value = 1802.8572
$ws.Range("I110").Value = 1847.6471
$ws.Range("J110").Value = 1612.5
$ws.Range("K110").Value = 1847.6471
$ws.Range("L110").Value = 1612.5
$ws.Range("M110").Value = 197.3529000000001
$ws.Range("N110").Value = -5702.5
$ws.Range("H122").Value = 2504
$ws.Range("I122").Value = 2006
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 6018
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -3568
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 2192.7693
$ws.Range("I132").Value = 2080.52
$ws.Range("K132").Value = 6241.559999999999
$ws.Range("M132").Value = -3711.559999999999
$ws.Range("H136").Value = 3158.7368
$ws.Range("I136").Value = 1768.3846
$ws.Range("J136").Value = 6171.1665
$ws.Range("K136").Value = 5305.1538
$ws.Range("L136").Value = 18513.4995
$ws.Range("M136").Value = -2755.1538
$ws.Range("N136").Value = -23613.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 24197.555
$ws.Range("J35").Value = 24197.555
$ws.Range("L35").Value = 24197.555
$ws.Range("N35").Value = -24817.555
$ws.Range("H82").Value = 47776.945
$ws.Range("I82").Value = 57759.727
$ws.Range("J82").Value = 32089.715
$ws.Range("K82").Value = 57759.727
$ws.Range("L82").Value = 32089.715
$ws.Range("M82").Value = -57376.727
$ws.Range("N82").Value = -32855.715
$ws.Range("H85").Value = 47776.945
$ws.Range("I85").Value = 57759.727
$ws.Range("J85").Value = 32089.715
$ws.Range("K85").Value = 57759.727
$ws.Range("L85").Value = 32089.715
$ws.Range("M85").Value = -56433.727
$ws.Range("N85").Value = -34741.715
$ws.Range("H99").Value = 3498.182
$ws.Range("I99").Value = 1650
$ws.Range("J99").Value = 4554.2856
$ws.Range("K99").Value = 1650
$ws.Range("L99").Value = 4554.2856
$ws.Range("M99").Value = -152
$ws.Range("N99").Value = -7550.2856
$ws.Range("H107").Value = 2460.3333
$ws.Range("I107").Value = 1960.579
$ws.Range("J107").Value = 4359.4
$ws.Range("K107").Value = 1960.579
$ws.Range("L107").Value = 4359.4
$ws.Range("M107").Value = -40.57899999999995
$ws.Range("N107").Value = -8199.4
$ws.Range("H134").Value = 2442.182
$ws.Range("I134").Value = 2353.7144
$ws.Range("K134").Value = 7061.1432
$ws.Range("M134").Value = -4526.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 15777.8
$ws.Range("J41").Value = 18232.5
$ws.Range("L41").Value = 18232.5
$ws.Range("N41").Value = -19088.5
$ws.Range("H51").Value = 18802.6
$ws.Range("J51").Value = 18802.6
$ws.Range("L51").Value = 18802.6
$ws.Range("N51").Value = -20274.6
$ws.Range("H60").Value = 15000
$ws.Range("J60").Value = 19500
$ws.Range("L60").Value = 19500
$ws.Range("N60").Value = -20522
$ws.Range("H61").Value = 18802.6
$ws.Range("J61").Value = 18802.6
$ws.Range("L61").Value = 18802.6
$ws.Range("N61").Value = -19498.6
$ws.Range("H134").Value = 6709.7144
$ws.Range("I134").Value = 7587.0557
$ws.Range("J134").Value = 1445.6666
$ws.Range("K134").Value = 22761.1671
$ws.Range("L134").Value = 4336.9998
$ws.Range("M134").Value = -20226.1671
$ws.Range("N134").Value = -9406.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 798.2
$ws.Range("I40").Value = 617.75
$ws.Range("J40").Value = 1520
$ws.Range("K40").Value = 2471
$ws.Range("L40").Value = 6080
$ws.Range("M40").Value = -2402
$ws.Range("N40").Value = -6218
$ws.Range("H68").Value = 2121.4783
$ws.Range("I68").Value = 1730.7778
$ws.Range("J68").Value = 2676.6843
$ws.Range("K68").Value = 5192.3334
$ws.Range("L68").Value = 8030.0529
$ws.Range("M68").Value = -4381.3334
$ws.Range("N68").Value = -9652.052899999999
$ws.Range("H71").Value = 2121.4783
$ws.Range("I71").Value = 1730.7778
$ws.Range("J71").Value = 2676.6843
$ws.Range("K71").Value = 15577.0002
$ws.Range("L71").Value = 24090.1587
$ws.Range("M71").Value = -11521.0002
$ws.Range("N71").Value = -32202.1587
$ws.Range("H107").Value = 846.2373
$ws.Range("I107").Value = 505.8095
$ws.Range("J107").Value = 1687.2941
$ws.Range("K107").Value = 1517.4285
$ws.Range("L107").Value = 5061.8823
$ws.Range("M107").Value = 402.5715
$ws.Range("N107").Value = -8901.882300000001
$ws.Range("H113").Value = 548.5
$ws.Range("I113").Value = 478.5
$ws.Range("J113").Value = 588.5
$ws.Range("K113").Value = 1435.5
$ws.Range("L113").Value = 1765.5
$ws.Range("M113").Value = 734.5
$ws.Range("N113").Value = -6105.5
$ws.Range("H117").Value = 3001.9333
$ws.Range("I117").Value = 405.8
$ws.Range("J117").Value = 4300
$ws.Range("K117").Value = 1217.4
$ws.Range("L117").Value = 12900
$ws.Range("M117").Value = 2224.6
$ws.Range("N117").Value = -19784

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14999.167
$ws.Range("I57").Value = 5238.5
$ws.Range("J57").Value = 27200
$ws.Range("K57").Value = 5238.5
$ws.Range("L57").Value = 27200
$ws.Range("M57").Value = -4418.5
$ws.Range("N57").Value = -28840
$ws.Range("H113").Value = 1744.4
$ws.Range("I113").Value = 1744.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1744.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 425.5999999999999
$ws.Range("N113").ClearContents()
$ws.Range("H123").Value = 16300
$ws.Range("J123").Value = 16300
$ws.Range("L123").Value = 16300
$ws.Range("N123").Value = -21200
$ws.Range("H135").Value = 49856.5
$ws.Range("I135").Value = 25000
$ws.Range("J135").Value = 52618.332
$ws.Range("K135").Value = 25000
$ws.Range("L135").Value = 52618.332
$ws.Range("M135").Value = -19930
$ws.Range("N135").Value = -62758.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2460.7
$ws.Range("I40").Value = 1601.1666
$ws.Range("J40").Value = 3750
$ws.Range("K40").Value = 1601.1666
$ws.Range("L40").Value = 3750
$ws.Range("M40").Value = -1465.1666
$ws.Range("N40").Value = -4022
$ws.Range("H136").Value = 1606.7727
$ws.Range("I136").Value = 1439.421
$ws.Range("J136").Value = 2666.6667
$ws.Range("K136").Value = 4318.263
$ws.Range("L136").Value = 8000.000100000001
$ws.Range("M136").Value = -1768.263
$ws.Range("N136").Value = -13100.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2792.9355
$ws.Range("I122").Value = 2499.6667
$ws.Range("J122").Value = 3408.8
$ws.Range("K122").Value = 7499.000100000001
$ws.Range("L122").Value = 10226.4
$ws.Range("M122").Value = -5049.000100000001
$ws.Range("N122").Value = -15126.4
$ws.Range("H132").Value = 3308.2046
$ws.Range("I132").Value = 3499.1
$ws.Range("J132").Value = 2899.1428
$ws.Range("K132").Value = 10497.3
$ws.Range("L132").Value = 8697.428400000001
$ws.Range("M132").Value = -7967.299999999999
$ws.Range("N132").Value = -13757.4284
